$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
$ws.Range("A4").Value = 104025305
$ws.Range("S4").Value = 10

$c = $ws.Range("Y4")
$c.NumberFormat = "@"
$c.Value = "2022-08-14"
$c.ClearFormats()

$ws.Range("Z4").Value = "11:44"

$c = $ws.Range("AA4")
$c.NumberFormat = "@"
$c.Value = "2022-08-14"
$c.ClearFormats()

$ws.Range("AB4").Value = "11:44"

$ws.Range("AO4").Value = "Sälg"
$ws.Range("AW4").Value = "Johan Staaf"
$ws.Range("AX4").Value = "Via Johan Staaf"
$ws.Range("AY4").Value = "LstZ inventering av skogliga värdetrakter 2022"

# --- Row 7 ---
$ws.Range("A7").Value = 104025268
$ws.Range("S7").Value = 10

$c = $ws.Range("Y7")
$c.NumberFormat = "@"
$c.Value = "2022-08-14"
$c.ClearFormats()

$ws.Range("Z7").Value = "11:51"

$c = $ws.Range("AA7")
$c.NumberFormat = "@"
$c.Value = "2022-08-14"
$c.ClearFormats()

$ws.Range("AB7").Value = "11:51"

$ws.Range("AW7").Value = "Johan Staaf"
$ws.Range("AX7").Value = "Via Johan Staaf"
$ws.Range("AY7").Value = "LstZ inventering av skogliga värdetrakter 2022"
